$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers - copy formatting (style) from H1, then set the header text
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows 2..33: I = 1 (constant), J = same as H
for ($r = 2; $r -le 33; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $hVal
}
